$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10, pushing the existing rows 10-87 down to 11-88.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly price record.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44761
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112021
$ws.Range("G10").Value = "Ají"
$ws.Range("H10").Value = "Inferno"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 18000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 18000
$ws.Range("N10").Value = "$/caja 15 kilos"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 1200
$ws.Range("Q10").Value = 15
$ws.Range("R10").Value = "Hortaliza"
